$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Week of 6/29/2018 (row 43): log Thursday's clock-in/out times
$ws.Range("C43").Value = 0.3
$ws.Range("F43").Value = 0.54375000000000007

# Week of 7/2/2018 (row 51): log Monday's clock-in/out times
$ws.Range("C51").Value = 0.25347222222222221
$ws.Range("F51").Value = 0.61111111111111105

# Grand total formula was pointing at the wrong weekly-total cell (G46, a
# blank spacer row) - fix it to reference the first week's total (G45)
$ws.Range("H57").Formula = "=SUM(G45, G57)"

# Reflect the scrolled viewport / active selection from the edit session
$ws.Range("H59").Select()
